$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Formulas")

# Row 3: drag the Sum/Average formulas from row 2 down into row 3
$ws.Range("H3").Formula = "=SUM(B3:G3)"
$ws.Range("I3").Formula = "=AVERAGE(B3:G3)"

# Row 4: drag the Profits formula (B4 = B2-B3) across columns C:G
$ws.Range("C4").Formula = "=C2-C3"
$ws.Range("D4").Formula = "=D2-D3"
$ws.Range("E4").Formula = "=E2-E3"
$ws.Range("F4").Formula = "=F2-F3"
$ws.Range("G4").Formula = "=G2-G3"

# Annotation cells describing what was dragged
$ws.Range("H4").Value = "Dragged to copy formula for profits"
$ws.Range("J2").Value = "Drag down column to copy avg and sum formula"

# Move the active selection like in the authored workbook
$ws.Range("C13").Select()
